$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newBiasVoltage = 0.846
$newRecombCurrent = [double]"7.0400000000000004E-5"

for ($r = 2; $r -le 37; $r++) {
    $ws.Cells.Item($r, 5).Value = $newBiasVoltage
    $ws.Cells.Item($r, 6).Value = $newRecombCurrent
}

$ws.Range("F2").Select()
